$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Add the new "Projected Resolution Date" column (E) to the Issues sheet.
$ws.Columns.Item(5).ColumnWidth = 24.84

$ws.Range("E1").Value = "Projected Resolution Date"

$ws.Range("E2").Value = 41564
$ws.Range("E2").NumberFormat = "mm-dd-yy"

# Make "Issues" the active / selected sheet, with C6 as the selected cell.
$ws.Activate()
$ws.Range("C6").Select()
